$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "Item Type" column (M) with a category per row.
# The six new category strings must be introduced (first write) in this exact
# order so the shared-strings table comes out in the same sequence as the
# authored workbook: Item Type, Breakfast, Entrée, Side, Desert, Drink.
# ---------------------------------------------------------------------------
$ws.Range("M1").Value  = "Item Type"
$ws.Range("M2").Value  = "Breakfast"
$ws.Range("M4").Value  = "Entrée"
$ws.Range("M5").Value  = "Side"
$ws.Range("M27").Value = "Desert"
$ws.Range("M21").Value = "Drink"

$categories = @{
    2  = "Breakfast"
    3  = "Breakfast"
    4  = "Entrée"
    5  = "Side"
    6  = "Entrée"
    7  = "Entrée"
    8  = "Entrée"
    9  = "Breakfast"
    10 = "Entrée"
    11 = "Entrée"
    12 = "Entrée"
    13 = "Side"
    14 = "Entrée"
    15 = "Side"
    16 = "Entrée"
    17 = "Entrée"
    18 = "Entrée"
    19 = "Entrée"
    20 = "Entrée"
    21 = "Drink"
    22 = "Breakfast"
    23 = "Drink"
    24 = "Breakfast"
    25 = "Drink"
    26 = "Drink"
    27 = "Desert"
    28 = "Desert"
    29 = "Drink"
    30 = "Drink"
    31 = "Entrée"
    32 = "Entrée"
    33 = "Entrée"
    34 = "Entrée"
    35 = "Entrée"
    36 = "Entrée"
    37 = "Side"
    38 = "Desert"
    39 = "Entrée"
    40 = "Entrée"
    41 = "Side"
    42 = "Side"
    43 = "Entrée"
    44 = "Entrée"
    45 = "Entrée"
    46 = "Side"
    47 = "Entrée"
    48 = "Desert"
    49 = "Side"
    50 = "Entrée"
}

for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 13).Value = $categories[$row]
}

# ---------------------------------------------------------------------------
# Dropdown list data validation on the whole column M.
# ---------------------------------------------------------------------------
$rng = $ws.Range("M1:M1048576")
$rng.Validation.Add(3, 1, 1, "Breakfast, Entrée, Desert, Drink, Side")
$rng.Validation.IgnoreBlank = $true
$rng.Validation.InCellDropdown = $true
$rng.Validation.ShowInput = $true
$rng.Validation.ShowError = $true

# ---------------------------------------------------------------------------
# Column widths (characters, Excel COM units) -- chosen so the serialized
# OOXML <col widths> come out matching the authored workbook.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth  = 38.666666666666664
$ws.Columns.Item(4).ColumnWidth  = 7.5
$ws.Columns.Item(5).ColumnWidth  = 8.666666666666666
$ws.Columns.Item(6).ColumnWidth  = 10.666666666666666
$ws.Columns.Item(7).ColumnWidth  = 10.333333333333334
$ws.Columns.Item(8).ColumnWidth  = 12.5
$ws.Columns.Item(9).ColumnWidth  = 12.833333333333334
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.666666666666666
$ws.Columns.Item(12).ColumnWidth = 7.0
$ws.Columns.Item(13).ColumnWidth = 8.0

# ---------------------------------------------------------------------------
# Sheet view: zoom to 150% and move the selection to O13.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 150
$ws.Range("O13").Select()
